$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.4258685247791129
$ws.Range("J2").Value = 0.5266590341921646
$ws.Range("M2").Value = 11.839004
$ws.Range("N2").Value = 35.517012
$ws.Range("O2").Value = 0.4816941403820247
$ws.Range("P2").Value = 0.5139419866672059
$ws.Range("Q2").Value = 4.090044389888001
$ws.Range("R2").Value = 36.810399508992
$ws.Range("S2").Value = 0.2051383729592358
$ws.Range("T2").Value = 0.270672190328953

# Row 3
$ws.Range("I3").Value = 0.4258685247791129
$ws.Range("J3").Value = 0.5266590341921646
$ws.Range("O3").Value = 0.3280382642169655
$ws.Range("P3").Value = 0.3499993524538634
$ws.Range("S3").Value = 0.1397011716531799
$ws.Range("T3").Value = 0.1843303209312347

# Row 4
$ws.Range("I4").Value = 0.4258685247791129
$ws.Range("J4").Value = 0.5266590341921646
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.049877
$ws.Range("N4").Value = 0.149631
$ws.Range("O4").Value = 0.002029347962027401
$ws.Range("P4").Value = 0.002165206166751885
$ws.Range("Q4").Value = 0.017231106944
$ws.Range("R4").Value = 0.155079962496
$ws.Range("S4").Value = 0.0008642354228521085
$ws.Range("T4").Value = 0.001140325388608467

# Row 5
$ws.Range("I5").Value = 0.4258685247791129
$ws.Range("J5").Value = 0.5266590341921646
$ws.Range("M5").Value = 4.6264905
$ws.Range("N5").Value = 9.252981
$ws.Range("O5").Value = 0.1882382474389825
$ws.Range("P5").Value = 0.1338934547121788
$ws.Range("Q5").Value = 1.598322926016
$ws.Range("R5").Value = 9.589937556096
$ws.Range("S5").Value = 0.08016474474384511
$ws.Range("T5").Value = 0.07051619754336839

# Row 6
$ws.Range("G6").Value = 0.4657455
$ws.Range("H6").Value = 0.931491
$ws.Range("I6").Value = 0.5741314752208871
$ws.Range("J6").Value = 0.4733409658078355
$ws.Range("M6").Value = 11.839004
$ws.Range("N6").Value = 35.517012
$ws.Range("O6").Value = 0.4816941403820247
$ws.Range("P6").Value = 0.5139419866672059
$ws.Range("Q6").Value = 5.513962837482
$ws.Range("R6").Value = 33.083777024892
$ws.Range("S6").Value = 0.276555767422789
$ws.Range("T6").Value = 0.2432697963382529

# Row 7
$ws.Range("G7").Value = 0.4657455
$ws.Range("H7").Value = 0.931491
$ws.Range("I7").Value = 0.5741314752208871
$ws.Range("J7").Value = 0.4733409658078355
$ws.Range("O7").Value = 0.3280382642169655
$ws.Range("P7").Value = 0.3499993524538634
$ws.Range("Q7").Value = 3.7550608291185
$ws.Range("R7").Value = 22.530364974711
$ws.Range("S7").Value = 0.1883370925637855
$ws.Range("T7").Value = 0.1656690315226287

# Row 8
$ws.Range("G8").Value = 0.4657455
$ws.Range("H8").Value = 0.931491
$ws.Range("I8").Value = 0.5741314752208871
$ws.Range("J8").Value = 0.4733409658078355
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.049877
$ws.Range("N8").Value = 0.149631
$ws.Range("O8").Value = 0.002029347962027401
$ws.Range("P8").Value = 0.002165206166751885
$ws.Range("Q8").Value = 0.0232299883035
$ws.Range("R8").Value = 0.139379929821
$ws.Range("S8").Value = 0.001165112539175293
$ws.Range("T8").Value = 0.001024880778143419

# Row 9
$ws.Range("G9").Value = 0.4657455
$ws.Range("H9").Value = 0.931491
$ws.Range("I9").Value = 0.5741314752208871
$ws.Range("J9").Value = 0.4733409658078355
$ws.Range("M9").Value = 4.6264905
$ws.Range("N9").Value = 9.252981
$ws.Range("O9").Value = 0.1882382474389825
$ws.Range("P9").Value = 0.1338934547121788
$ws.Range("Q9").Value = 2.15476713116775
$ws.Range("R9").Value = 8.619068524671
$ws.Range("S9").Value = 0.1080735026951374
$ws.Range("T9").Value = 0.06337725716881037
